# Fix the typo "gps" -> "GPS" in the subtitle of slide 1.
# The original text is split across three runs:
#   "The " + "gps" (err="1") + " that surrounds you"
# Using TextRange.Replace collapses the whole paragraph into a single run
# (matching how PowerPoint re-types/retypes a fully-selected text range),
# which is what the canonical OOXML after the edit looks like.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange
[void]$tr.Replace("The gps that surrounds you", "The GPS that surrounds you", 0, 0, 0)
